$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 change: proposal number
$ws.Range("A2").Value2 = "2363"

# Row 1 header changes: new "Numero Propuesta" / "Resultado" columns
$ws.Range("V1").Value2 = "Numero Propuesta"
$ws.Range("W1").Value2 = "Resultado"

# Row 2: result of the proposal, and drop the old Pass/Fail value
$ws.Range("W2").Value2 = "ok"
$ws.Range("V2").ClearContents()

# New Row 3 (duplicate of original row 2 content, but A3 = 2240)
$ws.Range("A3").Value2 = "2240"

# Copy formatting from row 2 (A:U) down to row 3, so text/number formats match
$ws.Range("A2:U2").Copy() | Out-Null
$ws.Range("A3:U3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in row 3 values (duplicate of row 2, but with the original proposal number)
$ws.Range("A3").Value2 = "2240"
$ws.Range("B3").Value2 = "CREDITO EMPRESARIAL"
$ws.Range("C3").Value2 = "CREDITOS PYMES"
$ws.Range("D3").Value2 = "NORMAL"
$ws.Range("E3").Value2 = "NORMAL"
$ws.Range("F3").Value2 = "SIN PROMOCION"
$ws.Range("G3").Value2 = "S/"
$ws.Range("H3").Value2 = "5000"
$ws.Range("I3").Value2 = "1"
$ws.Range("J3").Value2 = "Fija Vencida"
$ws.Range("K3").Value2 = "Libre Amortizacion"
$ws.Range("L3").Value2 = "90"
$ws.Range("M3").Value2 = "30"
$ws.Range("N3").Value2 = "GIRO BANCO DE LA NACION"
$ws.Range("O3").Value2 = "LIMA"
$ws.Range("P3").Value2 = "LIMA"
$ws.Range("Q3").Value2 = "LIMA"
$ws.Range("R3").Value2 = "prueba de nueva propuesta"
$ws.Range("S3").Value2 = "prueba de nueva propuesta"
$ws.Range("T3").Value2 = "Aprobar propuesta"
$ws.Range("U3").Value2 = "prueba"
$ws.Range("W3").Value2 = "ok"

$ws.Range("B7").Select() | Out-Null
